$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("transitVehicleToCapacity")
$ws.Activate()

# Insert two new rows at row 4, pushing the existing "7 Car BART" block (and
# everything below it) down by two rows.
$ws.Range("A4:A5").EntireRow.Insert()

# Row 4: 5 Car BART
$ws.Cells.Item(4,1).Value = "5 Car BART"
$ws.Cells.Item(4,2).Value = 555
$ws.Cells.Item(4,3).Value = 471.75
$ws.Cells.Item(4,4).Value = "5 Car BART"
$ws.Cells.Item(4,5).Value = 0
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 0
$ws.Cells.Item(4,8).Value = 0

# Row 5: 5 Car BART RENOVATED
$ws.Cells.Item(5,1).Value = "5 Car BART RENOVATED"
$ws.Cells.Item(5,2).Value = 560
$ws.Cells.Item(5,3).Value = 476
$ws.Cells.Item(5,4).Value = "5 Car BART RENOVATED"
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = 0
$ws.Cells.Item(5,7).Value = 0
$ws.Cells.Item(5,8).Value = 0

# Match the author's final view/selection state on the sheet.
$ws.Range("D7").Select()
